$d = $word.ActiveDocument

# The test-plan table is the only table in the document; grab it and
# append a brand new row for Test 10 ("Test GetLength and GetSize"),
# matching the style of every other row already in the table.
$t = $d.Tables.Item(1)
$newRow = $t.Rows.Add()
$rowIndex = $newRow.Index

# Column 1 - Test number
$cell1 = $t.Cell($rowIndex, 1)
$cell1.Range.Text = "10"

# Column 2 - Description
$cell2 = $t.Cell($rowIndex, 2)
$cell2.Range.Text = "Test GetLength and GetSize"

# Column 3 - Change code (several lines, including two blank ones)
$cell3 = $t.Cell($rowIndex, 3)
$code = "Vector<int> testVec(3);`r" + `
        "`r" + `
        "for(int i = 0; i < 2; i++)`r" + `
        "     testVec.PushBack(i);`r" + `
        "`r" + `
        "cout << testVec.GetLength() << endl;`r" + `
        "cout << testVec.GetSize() << endl;"
$cell3.Range.Text = $code

# Column 4 - Expected output
$cell4 = $t.Cell($rowIndex, 4)
$cell4.Range.Text = "2`r3"

# Column 5 - Passed
$cell5 = $t.Cell($rowIndex, 5)
$cell5.Range.Text = "Y"
